# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de and
# zh-cn handback packages have now come back from localization:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#    columns on the per-language sheets get populated
#  - A few columns are widened so the new, longer values aren't clipped

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$targetFileName  = "1124f0ba-16a7-43c8-8907-272dd276cb4a.md"
$targetUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9344e470d03b0a94fb3a185125c86128d897608f/e2e/1124f0ba-16a7-43c8-8907-272dd276cb4a.md"

$zhHandbackFile  = "1124f0ba-16a7-43c8-8907-272dd276cb4a.b00ae37202c4d966e8b28155bf0cca554205fcf0.zh-cn.xlf"
$deHandbackFile  = "1124f0ba-16a7-43c8-8907-272dd276cb4a.b00ae37202c4d966e8b28155bf0cca554205fcf0.de-de.xlf"

$zhHandbackDate  = "2016-08-19 15:09:29"
$deHandbackDate  = "2016-08-19 15:09:36"

# Column width helper: the widest column used elsewhere in this workbook is
# rendered as XML width="40", which corresponds to a COM ColumnWidth of
# 39.166666... in this engine. The new "wide" columns (29.9777... in the
# target) round to the same quantized column-width step as ColumnWidth
# 29.166666..., so we reuse these two constants everywhere below.
$wideWidth   = 39.166666666666664
$mediumWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value2 = $statusNew
$wsOverview.Range("F2").Value2 = $statusNew
$wsOverview.Range("E3").Value2 = $statusNew
$wsOverview.Range("F3").Value2 = $statusNew

$wsOverview.Columns.Item(5).ColumnWidth = $mediumWidth
$wsOverview.Columns.Item(6).ColumnWidth = $mediumWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value2 = $statusNew
$wsZhCn.Range("C3").Value2 = $statusNew

$wsZhCn.Columns.Item(3).ColumnWidth  = $mediumWidth
$wsZhCn.Columns.Item(9).ColumnWidth  = $wideWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideWidth

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName)

$wsZhCn.Range("J2").Value2 = $zhHandbackFile
$wsZhCn.Range("J3").Value2 = $zhHandbackFile

$wsZhCn.Range("K2").Value2 = $zhHandbackDate
$wsZhCn.Range("K3").Value2 = $zhHandbackDate

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value2 = $statusNew
$wsDeDe.Range("C3").Value2 = $statusNew

$wsDeDe.Columns.Item(3).ColumnWidth  = $mediumWidth
$wsDeDe.Columns.Item(9).ColumnWidth  = $wideWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideWidth

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName)

$wsDeDe.Range("J2").Value2 = $deHandbackFile
$wsDeDe.Range("J3").Value2 = $deHandbackFile

$wsDeDe.Range("K2").Value2 = $deHandbackDate
$wsDeDe.Range("K3").Value2 = $deHandbackDate
